# Regenerate merged AHB files
# - rename the "_old"/"_new" header suffixes to "_FV2410"/"_FV2504"
# - freeze the header row
# - turn the data range into a real Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) rename header row cells (columns A:J were "..._old", L:U were "..._new") ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("K1").Value = "diff"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2) freeze the header row (split below row 1, freeze panes) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3) convert the used range into a native Excel table ---
$rng = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
